$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 109.026058
$ws.Range("H2").Value = 327.078174
$ws.Range("I2").Value = 0.3049840938689738
$ws.Range("J2").Value = 0.3049840938689738
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.214110666666667
$ws.Range("N2").Value = 21.642332
$ws.Range("O2").Value = 0.4688823795981188
$ws.Range("P2").Value = 0.4688823795981188
$ws.Range("Q2").Value = 786.5260479624186
$ws.Range("R2").Value = 7078.734431661767
$ws.Range("S2").Value = 0.1430016676728605
$ws.Range("T2").Value = 0.1430016676728605

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 109.026058
$ws.Range("H3").Value = 327.078174
$ws.Range("I3").Value = 0.3049840938689738
$ws.Range("J3").Value = 0.3049840938689738
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.110350666666666
$ws.Range("N3").Value = 21.331052
$ws.Range("O3").Value = 0.4621384803214003
$ws.Range("P3").Value = 0.4621384803214003
$ws.Range("Q3").Value = 775.2135041843386
$ws.Range("R3").Value = 6976.921537659047
$ws.Range("S3").Value = 0.1409448856628069
$ws.Range("T3").Value = 0.1409448856628069

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 109.026058
$ws.Range("H4").Value = 327.078174
$ws.Range("I4").Value = 0.3049840938689738
$ws.Range("J4").Value = 0.3049840938689738
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.061296333333333
$ws.Range("N4").Value = 3.183889
$ws.Range("O4").Value = 0.06897914008048092
$ws.Range("P4").Value = 0.06897914008048092
$ws.Range("Q4").Value = 115.7089555931873
$ws.Range("R4").Value = 1041.380600338686
$ws.Range("S4").Value = 0.02103754053330649
$ws.Range("T4").Value = 0.02103754053330649

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 89.97721833333333
$ws.Range("H5").Value = 269.931655
$ws.Range("I5").Value = 0.2516978134001918
$ws.Range("J5").Value = 0.2516978134001917
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.214110666666667
$ws.Range("N5").Value = 21.642332
$ws.Range("O5").Value = 0.4688823795981188
$ws.Range("P5").Value = 0.4688823795981188
$ws.Range("Q5").Value = 649.1056105354955
$ws.Range("R5").Value = 5841.950494819459
$ws.Range("S5").Value = 0.1180166696867252
$ws.Range("T5").Value = 0.1180166696867252

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 89.97721833333333
$ws.Range("H6").Value = 269.931655
$ws.Range("I6").Value = 0.2516978134001918
$ws.Range("J6").Value = 0.2516978134001917
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.110350666666666
$ws.Range("N6").Value = 21.331052
$ws.Range("O6").Value = 0.4621384803214003
$ws.Range("P6").Value = 0.4621384803214003
$ws.Range("Q6").Value = 639.7695743612288
$ws.Range("R6").Value = 5757.926169251059
$ws.Range("S6").Value = 0.116319244984984
$ws.Range("T6").Value = 0.116319244984984

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 89.97721833333333
$ws.Range("H7").Value = 269.931655
$ws.Range("I7").Value = 0.2516978134001918
$ws.Range("J7").Value = 0.2516978134001917
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.061296333333333
$ws.Range("N7").Value = 3.183889
$ws.Range("O7").Value = 0.06897914008048092
$ws.Range("P7").Value = 0.06897914008048092
$ws.Range("Q7").Value = 95.49249190069942
$ws.Range("R7").Value = 859.4324271062949
$ws.Range("S7").Value = 0.01736189872848258
$ws.Range("T7").Value = 0.01736189872848257

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 158.477852
$ws.Range("H8").Value = 475.433556
$ws.Range("I8").Value = 0.4433180927308344
$ws.Range("J8").Value = 0.4433180927308344
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.214110666666667
$ws.Range("N8").Value = 21.642332
$ws.Range("O8").Value = 0.4688823795981188
$ws.Range("P8").Value = 0.4688823795981188
$ws.Range("Q8").Value = 1143.276762543621
$ws.Range("R8").Value = 10289.49086289259
$ws.Range("S8").Value = 0.2078640422385332
$ws.Range("T8").Value = 0.2078640422385331

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 158.477852
$ws.Range("H9").Value = 475.433556
$ws.Range("I9").Value = 0.4433180927308344
$ws.Range("J9").Value = 0.4433180927308344
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.110350666666666
$ws.Range("N9").Value = 21.331052
$ws.Range("O9").Value = 0.4621384803214003
$ws.Range("P9").Value = 0.4621384803214003
$ws.Range("Q9").Value = 1126.833100620101
$ws.Range("R9").Value = 10141.49790558091
$ws.Range("S9").Value = 0.2048743496736095
$ws.Range("T9").Value = 0.2048743496736094

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 158.477852
$ws.Range("H10").Value = 475.433556
$ws.Range("I10").Value = 0.4433180927308344
$ws.Range("J10").Value = 0.4433180927308344
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.061296333333333
$ws.Range("N10").Value = 3.183889
$ws.Range("O10").Value = 0.06897914008048092
$ws.Range("P10").Value = 0.06897914008048092
$ws.Range("Q10").Value = 168.1919632421426
$ws.Range("R10").Value = 1513.727669179284
$ws.Range("S10").Value = 0.03057970081869186
$ws.Range("T10").Value = 0.03057970081869186

